# Insert a new weekly price record as row 160 in the "Femacal de La Calera - Ciboulette"
# sheet, shifting the existing rows 160:264 down to 161:265.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 160, shifting rows 160:264 -> 161:265.
$ws.Rows(160).Insert()

# Populate the newly inserted row 160 with the new record.
$ws.Range("A160").Value2 = 3
$ws.Range("B160").Value2 = "Femacal de La Calera"
$ws.Range("C160").Value2 = "Coquimbo"
$ws.Range("D160").Value2 = 44603
$ws.Range("E160").Value2 = 5
$ws.Range("F160").Value2 = 100112039
$ws.Range("G160").Value2 = "Ciboulette"
$ws.Range("H160").Value2 = "Sin especificar"
$ws.Range("I160").Value2 = "Primera"
$ws.Range("J160").Value2 = 160
$ws.Range("K160").Value2 = 1500
$ws.Range("L160").Value2 = 1500
$ws.Range("M160").Value2 = 1500
$ws.Range("N160").Value2 = "$/docena de atados"
$ws.Range("O160").Value2 = "Provincia de Quillota"
$ws.Range("P160").Value2 = 500
$ws.Range("Q160").Value2 = 3
$ws.Range("R160").Value2 = "Hortaliza"
